$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.951.11"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "3.443.27"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'583.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").Value = "'173.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "3.441.47"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").Value = "'6.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").Value = "4.041.62"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "'28.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.35%  "
$ws.Range("D16").Value = "66.022.82"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "3.441.46"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "'5.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("D20").Value = "'13.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").Value = "'368.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("D22").Value = "'7.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "'72.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'0.535"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "'9.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'23.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.07%  "
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").Value = "'5.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.18%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "'1.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.38%  "
$ws.Range("D35").Value = "'7.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").Value = "'160.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'28.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("D39").Value = "'0.880"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("D42").Value = "2.756.85"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'6.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "'0.0680"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").Value = "'40.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").Value = "'24.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D49").Value = "'324.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "'6.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  -2.59%  "

Write-Host "Applied 96 cell updates"
